$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, shifting the existing rows 7..41 down to 8..42
# (mirrors Excel's "Insert Sheet Rows" which pushes data down and carries the
# formatting of the row above, which is why column D keeps its date style).
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new weekly price observation.
$ws.Cells.Item(7, 1).Value = 7
$ws.Cells.Item(7, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(7, 3).Value = "Ñuble"
$ws.Cells.Item(7, 4).Value = 45022
$ws.Cells.Item(7, 5).Value = 16
$ws.Cells.Item(7, 6).Value = 100112044
$ws.Cells.Item(7, 7).Value = "Perejil"
$ws.Cells.Item(7, 8).Value = "Sin especificar"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 230
$ws.Cells.Item(7, 11).Value = 1400
$ws.Cells.Item(7, 12).Value = 1500
$ws.Cells.Item(7, 13).Value = 1465
$ws.Cells.Item(7, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(7, 15).Value = "Región del Maule"
$ws.Cells.Item(7, 16).Value = 1465
$ws.Cells.Item(7, 17).Value = 1
$ws.Cells.Item(7, 18).Value = "Hortaliza"
